# Auto-generated edit script applying the Zeromus_Profits market-data refresh diff
# across the ALC/ARM/CRP/CUL/LTW/WVR sheets (leve profit calculator).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)  # ALC
$ws.Cells.Item(3,8).Value = 29978
$ws.Cells.Item(3,10).Value = 29978
$ws.Cells.Item(3,12).Value = 29978
$ws.Cells.Item(3,14).Value = -30206

$ws.Cells.Item(31,8).Value = 300.5
$ws.Cells.Item(31,9).Value = 300.5
$ws.Cells.Item(31,11).Value = 901.5
$ws.Cells.Item(31,13).Value = -671.5

$ws.Cells.Item(102,8).Value = 29978
$ws.Cells.Item(102,10).Value = 29978
$ws.Cells.Item(102,12).Value = 29978
$ws.Cells.Item(102,14).Value = -36468

$ws.Cells.Item(127,8).Value = 802.5714
$ws.Cells.Item(127,9).Value = 603.5
$ws.Cells.Item(127,10).Value = 1997
$ws.Cells.Item(127,11).Value = 1810.5
$ws.Cells.Item(127,12).Value = 5991
$ws.Cells.Item(127,13).Value = 3149.5
$ws.Cells.Item(127,14).Value = -15911

$ws.Cells.Item(129,8).Value = 8083.033
$ws.Cells.Item(129,9).Value = 397.7
$ws.Cells.Item(129,10).Value = 11925.7
$ws.Cells.Item(129,11).Value = 1193.1
$ws.Cells.Item(129,12).Value = 35777.10000000001
$ws.Cells.Item(129,13).Value = 3806.9
$ws.Cells.Item(129,14).Value = -45777.10000000001

$ws.Cells.Item(132,8).Value = 5667.381
$ws.Cells.Item(132,9).Value = 5919.1177
$ws.Cells.Item(132,10).Value = 4597.5
$ws.Cells.Item(132,11).Value = 17757.3531
$ws.Cells.Item(132,12).Value = 13792.5
$ws.Cells.Item(132,13).Value = -15227.3531
$ws.Cells.Item(132,14).Value = -18852.5

$ws.Cells.Item(138,8).Value = 3414.3552
$ws.Cells.Item(138,9).Value = 2660.2354
$ws.Cells.Item(138,10).Value = 3631.644
$ws.Cells.Item(138,11).Value = 7980.706200000001
$ws.Cells.Item(138,12).Value = 10894.932
$ws.Cells.Item(138,13).Value = -2840.706200000001
$ws.Cells.Item(138,14).Value = -21174.932

$ws = $wb.Worksheets.Item(2)  # ARM
$ws.Cells.Item(61,8).Value = 2278.8635
$ws.Cells.Item(61,9).Value = 2243.5
$ws.Cells.Item(61,11).Value = 2243.5
$ws.Cells.Item(61,13).Value = -2031.5

$ws.Cells.Item(92,8).Value = 29275
$ws.Cells.Item(92,10).Value = 29275
$ws.Cells.Item(92,12).Value = 29275
$ws.Cells.Item(92,14).Value = -34267

$ws.Cells.Item(95,8).Value = 22200
$ws.Cells.Item(95,10).Value = 22200
$ws.Cells.Item(95,12).Value = 22200
$ws.Cells.Item(95,14).Value = -27692

$ws.Cells.Item(96,8).Value = 27562.666
$ws.Cells.Item(96,10).Value = 27562.666
$ws.Cells.Item(96,12).Value = 27562.666
$ws.Cells.Item(96,14).Value = -33054.666

$ws.Cells.Item(109,8).Value = 28725.666
$ws.Cells.Item(109,10).Value = 28725.666
$ws.Cells.Item(109,12).Value = 28725.666
$ws.Cells.Item(109,14).Value = -31499.666

$ws.Cells.Item(110,8).Value = 1426.5758
$ws.Cells.Item(110,9).Value = 1411.2916
$ws.Cells.Item(110,11).Value = 1411.2916
$ws.Cells.Item(110,13).Value = 633.7084

$ws.Cells.Item(132,8).Value = 41008.668
$ws.Cells.Item(132,9).Value = 49356.184
$ws.Cells.Item(132,10).Value = 4279.6
$ws.Cells.Item(132,11).Value = 148068.552
$ws.Cells.Item(132,12).Value = 12838.8
$ws.Cells.Item(132,13).Value = -145538.552
$ws.Cells.Item(132,14).Value = -17898.8

$ws.Cells.Item(136,8).Value = 2278.8635
$ws.Cells.Item(136,9).Value = 2243.5
$ws.Cells.Item(136,11).Value = 6730.5
$ws.Cells.Item(136,13).Value = -4180.5

$ws = $wb.Worksheets.Item(4)  # CRP
$ws.Cells.Item(23,8).Value = 1100
$ws.Cells.Item(23,9).Value = 1100
$ws.Cells.Item(23,10).Value = 0
$ws.Cells.Item(23,11).Value = 1100
$ws.Cells.Item(23,12).Value = 0
$ws.Cells.Item(23,13).Value = -860
$ws.Cells.Item(23,14).ClearContents()

$ws.Cells.Item(27,8).Value = 1100
$ws.Cells.Item(27,9).Value = 1100
$ws.Cells.Item(27,10).Value = 0
$ws.Cells.Item(27,11).Value = 1100
$ws.Cells.Item(27,12).Value = 0
$ws.Cells.Item(27,13).Value = -908
$ws.Cells.Item(27,14).ClearContents()

$ws.Cells.Item(31,8).Value = 2019.3115
$ws.Cells.Item(31,9).Value = 1113.9286
$ws.Cells.Item(31,10).Value = 2787.5151
$ws.Cells.Item(31,11).Value = 1113.9286
$ws.Cells.Item(31,12).Value = 2787.5151
$ws.Cells.Item(31,13).Value = -818.9286
$ws.Cells.Item(31,14).Value = -3377.5151

$ws.Cells.Item(34,8).Value = 2019.3115
$ws.Cells.Item(34,9).Value = 1113.9286
$ws.Cells.Item(34,10).Value = 2787.5151
$ws.Cells.Item(34,11).Value = 1113.9286
$ws.Cells.Item(34,12).Value = 2787.5151
$ws.Cells.Item(34,13).Value = -911.9286
$ws.Cells.Item(34,14).Value = -3191.5151

$ws = $wb.Worksheets.Item(5)  # CUL
$ws.Cells.Item(5,8).Value = 816.6389
$ws.Cells.Item(5,9).Value = 506.30435
$ws.Cells.Item(5,11).Value = 1518.91305
$ws.Cells.Item(5,13).Value = -1406.91305

$ws.Cells.Item(20,8).Value = 1590
$ws.Cells.Item(20,10).Value = 2000
$ws.Cells.Item(20,12).Value = 6000
$ws.Cells.Item(20,14).Value = -6454

$ws.Cells.Item(24,8).Value = 543.3333
$ws.Cells.Item(24,10).Value = 690
$ws.Cells.Item(24,12).Value = 2070
$ws.Cells.Item(24,14).Value = -2530

$ws.Cells.Item(35,8).Value = 0
$ws.Cells.Item(35,9).Value = 0
$ws.Cells.Item(35,11).Value = 0
$ws.Cells.Item(35,13).ClearContents()

$ws.Cells.Item(37,8).Value = 59086.668
$ws.Cells.Item(37,10).Value = 59086.668
$ws.Cells.Item(37,12).Value = 177260.004
$ws.Cells.Item(37,14).Value = -177484.004

$ws.Cells.Item(131,8).Value = 830.84906
$ws.Cells.Item(131,9).Value = 499.05554
$ws.Cells.Item(131,10).Value = 1001.4857
$ws.Cells.Item(131,11).Value = 1497.16662
$ws.Cells.Item(131,12).Value = 3004.4571
$ws.Cells.Item(131,13).Value = 3542.83338
$ws.Cells.Item(131,14).Value = -13084.4571

$ws.Cells.Item(132,8).Value = 1099
$ws.Cells.Item(132,9).Value = 527.6667
$ws.Cells.Item(132,10).Value = 1813.1666
$ws.Cells.Item(132,11).Value = 4749.0003
$ws.Cells.Item(132,12).Value = 16318.4994
$ws.Cells.Item(132,13).Value = -2219.0003
$ws.Cells.Item(132,14).Value = -21378.4994

$ws.Cells.Item(135,8).Value = 816.6389
$ws.Cells.Item(135,9).Value = 506.30435
$ws.Cells.Item(135,11).Value = 4556.73915
$ws.Cells.Item(135,13).Value = -2021.73915

$ws = $wb.Worksheets.Item(7)  # LTW
$ws.Cells.Item(16,8).Value = 1123.75
$ws.Cells.Item(16,9).Value = 1123.75
$ws.Cells.Item(16,11).Value = 1123.75
$ws.Cells.Item(16,13).Value = -953.75

$ws = $wb.Worksheets.Item(8)  # WVR
$ws.Cells.Item(61,8).Value = 6202.8
$ws.Cells.Item(61,9).Value = 1000
$ws.Cells.Item(61,10).Value = 7503.5
$ws.Cells.Item(61,11).Value = 1000
$ws.Cells.Item(61,12).Value = 7503.5
$ws.Cells.Item(61,13).Value = -708
$ws.Cells.Item(61,14).Value = -8087.5

$ws.Cells.Item(63,8).Value = 26000
$ws.Cells.Item(63,10).Value = 26000
$ws.Cells.Item(63,12).Value = 26000
$ws.Cells.Item(63,14).Value = -27248

$ws.Cells.Item(66,8).Value = 26000
$ws.Cells.Item(66,10).Value = 26000
$ws.Cells.Item(66,12).Value = 78000
$ws.Cells.Item(66,14).Value = -84240

$ws.Cells.Item(80,8).Value = 28850
$ws.Cells.Item(80,10).Value = 28850
$ws.Cells.Item(80,12).Value = 28850
$ws.Cells.Item(80,14).Value = -30846

$ws.Cells.Item(83,8).Value = 28850
$ws.Cells.Item(83,10).Value = 28850
$ws.Cells.Item(83,12).Value = 86550
$ws.Cells.Item(83,14).Value = -96534

$ws.Cells.Item(98,8).Value = 28500
$ws.Cells.Item(98,10).Value = 28500
$ws.Cells.Item(98,12).Value = 28500
$ws.Cells.Item(98,14).Value = -34490

$ws.Cells.Item(132,8).Value = 4134.5
$ws.Cells.Item(132,9).Value = 1601.3334
$ws.Cells.Item(132,10).Value = 6667.6665
$ws.Cells.Item(132,11).Value = 4804.0002
$ws.Cells.Item(132,12).Value = 20002.9995
$ws.Cells.Item(132,13).Value = -2274.0002
$ws.Cells.Item(132,14).Value = -25062.9995
